$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---- Row 10: Tree / Hard / 297. Serialize and Deserialize Binary Tree ----
$ws.Cells.Item(10,1).Value = "Tree"
$ws.Cells.Item(10,2).Value = "Hard"
$name10 = "297. Serialize and Deserialize Binary Tree"
$ws.Cells.Item(10,3).Value = $name10
$ws.Cells.Item(10,4).Value = 'First serialize into an array using preorder traversal and then return it as string with ",".join(serialized). Next split the serialized over "," and iterate over values (maybe store the vals into an iterator with vals = iter(data) and rebuild the same way you serialized using preorder'

$url10 = "https://leetcode.com/problems/serialize-and-deserialize-binary-tree/"
$h10 = $ws.Hyperlinks.Add($ws.Range("C10"), $url10)
$h10.TextToDisplay = $url10
$ws.Cells.Item(10,3).Value = $name10
$ws.Range("C10").Style = "Neutral"

$ws.Rows.Item(10).RowHeight = 33

# ---- Row 11: Heap/PQ / Easy / 703. Kth Largest Element in a Stream ----
$ws.Cells.Item(11,1).Value = "Heap/PQ"
$ws.Cells.Item(11,2).Value = "Easy"
$name11 = "703. Kth Largest Element in a Stream"
$ws.Cells.Item(11,3).Value = $name11
$ws.Cells.Item(11,4).Value = "Use a heap to track top k elements (remove any elements after k for optmization)"

$url11 = "https://leetcode.com/problems/kth-largest-element-in-a-stream/"
$h11 = $ws.Hyperlinks.Add($ws.Range("C11"), $url11)
$h11.TextToDisplay = $url11
$ws.Cells.Item(11,3).Value = $name11
$ws.Range("C11").Style = "Good"

# ---- Row 12: Heap/PQ / Easy / 1046. Last Stone Weight ----
$ws.Cells.Item(12,1).Value = "Heap/PQ"
$ws.Cells.Item(12,2).Value = "Easy"
$name12 = "1046. Last Stone Weight"
$ws.Cells.Item(12,3).Value = $name12
$ws.Cells.Item(12,4).Value = "Just use a heap. Re add elements to the heap again after using them (if needed). Be careful with the negatives since python max heap needs you to flip the signs"

$url12 = "https://leetcode.com/problems/last-stone-weight/"
$h12 = $ws.Hyperlinks.Add($ws.Range("C12"), $url12)
$h12.TextToDisplay = $url12
$ws.Cells.Item(12,3).Value = $name12
$ws.Range("C12").Style = "Good"

$ws.Rows.Item(12).RowHeight = 28.8

# ---- View / selection state to match the saved workbook ----
$ws.Select()
$ws.Range("D12").Select()
$excel.ActiveWindow.ScrollRow = 8
$excel.ActiveWindow.ScrollColumn = 1

Write-Host "done"
